$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values must remain stored as text, matching the source inlineStr cells.
# Set each cell's NumberFormat to text, assign the literal string, then restore the default
# (un-styled) cell style so no stray style index is left behind.
$numericTextValues = @{
    'D2' = '279.17'
    'E2' = '0.43%'
    'D3' = '27.45'
    'E3' = '0.36%'
    'D4' = '4.834'
    'E4' = '-1.36%'
    'D5' = '0.06381'
    'E5' = '0.14%'
    'D6' = '7.023'
    'E6' = '0.68%'
    'D7' = '1.306'
    'E7' = '3.19%'
    'D8' = '0.8939'
    'E8' = '1.27%'
    'D9' = '0.1525'
    'E9' = '-0.15%'
    'D10' = '0.05894'
    'E10' = '15.37%'
    'E11' = '-0.44%'
    'D12' = '0.02909'
    'E12' = '-1.89%'
    'D13' = '0.08999'
    'E13' = '-0.09%'
    'D14' = '0.001575'
    'E14' = '0.52%'
    'D15' = '0.0006416'
    'E15' = '0.09%'
    'D16' = '0.006107'
    'E16' = '2.38%'
    'E17' = '0.50%'
    'D18' = '3.326'
    'E18' = '0.45%'
    'D19' = '2.226'
    'E19' = '-2.02%'
    'E21' = '1.11%'
    'D22' = '3.898'
    'E22' = '-0.49%'
    'D23' = '0.04411'
    'E23' = '-0.24%'
    'D24' = '0.1504'
    'E24' = '8.95%'
    'E25' = '0.33%'
    'D26' = '0.004276'
    'E26' = '10.45%'
    'E28' = '-1.67%'
    'D29' = '0.0001654'
    'E29' = '-14.58%'
    'D40' = '0.04074'
    'E40' = '-1.75%'
    'D41' = '0.006646'
    'E41' = '-3.02%'
    'D42' = '0.1398'
    'E42' = '18.64%'
    'D43' = '0.002098'
    'E43' = '3.96%'
    'E44' = '-2.00%'
    'D45' = '0.00005559'
    'E45' = '7.61%'
    'E46' = '5.00%'
    'D47' = '0.01850'
    'E47' = '-8.57%'
}
foreach ($cellRef in $numericTextValues.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextValues[$cellRef]
    $cell.Style = "Normal"
}

# Plain text fields (coin name / link) can be set directly
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
